$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "H2-M3"
$ws.Cells.Item(2, 3).Value = "Klrd1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.249240666666667
$ws.Cells.Item(2, 8).Value = 6.747722
$ws.Cells.Item(2, 9).Value = 0.1370621209892393
$ws.Cells.Item(2, 10).Value = 0.1370621209892393
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 1.286367333333333
$ws.Cells.Item(2, 14).Value = 3.859102
$ws.Cells.Item(2, 15).Value = 0.416990774147385
$ws.Cells.Item(2, 16).Value = 0.4169907741473849
$ws.Cells.Item(2, 17).Value = 2.893349718404889
$ws.Cells.Item(2, 18).Value = 26.040147465644
$ws.Cells.Item(2, 19).Value = 0.05715363993758545
$ws.Cells.Item(2, 20).Value = 0.05715363993758543

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "H2-M3"
$ws.Cells.Item(3, 3).Value = "Klrd1"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.249240666666667
$ws.Cells.Item(3, 8).Value = 6.747722
$ws.Cells.Item(3, 9).Value = 0.1370621209892393
$ws.Cells.Item(3, 10).Value = 0.1370621209892393
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.798514666666667
$ws.Cells.Item(3, 14).Value = 5.395544
$ws.Cells.Item(3, 15).Value = 0.5830092258526149
$ws.Cells.Item(3, 16).Value = 0.5830092258526149
$ws.Cells.Item(3, 17).Value = 4.045292327863111
$ws.Cells.Item(3, 18).Value = 36.407630950768
$ws.Cells.Item(3, 19).Value = 0.07990848105165385
$ws.Cells.Item(3, 20).Value = 0.07990848105165384

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "H2-M3"
$ws.Cells.Item(4, 3).Value = "Klrd1"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.061273333333333
$ws.Cells.Item(4, 8).Value = 6.18382
$ws.Cells.Item(4, 9).Value = 0.1256079436905785
$ws.Cells.Item(4, 10).Value = 0.1256079436905785
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 1.286367333333333
$ws.Cells.Item(4, 14).Value = 3.859102
$ws.Cells.Item(4, 15).Value = 0.416990774147385
$ws.Cells.Item(4, 16).Value = 0.4169907741473849
$ws.Cells.Item(4, 17).Value = 2.651554681071111
$ws.Cells.Item(4, 18).Value = 23.86399212964
$ws.Cells.Item(4, 19).Value = 0.05237735367859548
$ws.Cells.Item(4, 20).Value = 0.05237735367859546

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "H2-M3"
$ws.Cells.Item(5, 3).Value = "Klrd1"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.061273333333333
$ws.Cells.Item(5, 8).Value = 6.18382
$ws.Cells.Item(5, 9).Value = 0.1256079436905785
$ws.Cells.Item(5, 10).Value = 0.1256079436905785
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.798514666666667
$ws.Cells.Item(5, 14).Value = 5.395544
$ws.Cells.Item(5, 15).Value = 0.5830092258526149
$ws.Cells.Item(5, 16).Value = 0.5830092258526149
$ws.Cells.Item(5, 17).Value = 3.707230322008889
$ws.Cells.Item(5, 18).Value = 33.36507289808
$ws.Cells.Item(5, 19).Value = 0.07323059001198302
$ws.Cells.Item(5, 20).Value = 0.07323059001198301

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "H2-M3"
$ws.Cells.Item(6, 3).Value = "Klrd1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.47742133333333
$ws.Cells.Item(6, 8).Value = 34.432264
$ws.Cells.Item(6, 9).Value = 0.6994003508593609
$ws.Cells.Item(6, 10).Value = 0.6994003508593608
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 1.286367333333333
$ws.Cells.Item(6, 14).Value = 3.859102
$ws.Cells.Item(6, 15).Value = 0.416990774147385
$ws.Cells.Item(6, 16).Value = 0.4169907741473849
$ws.Cells.Item(6, 17).Value = 14.76417987410311
$ws.Cells.Item(6, 18).Value = 132.877618866928
$ws.Cells.Item(6, 19).Value = 0.2916434937437976
$ws.Cells.Item(6, 20).Value = 0.2916434937437975

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "H2-M3"
$ws.Cells.Item(7, 3).Value = "Klrd1"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.47742133333333
$ws.Cells.Item(7, 8).Value = 34.432264
$ws.Cells.Item(7, 9).Value = 0.6994003508593609
$ws.Cells.Item(7, 10).Value = 0.6994003508593608
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.798514666666667
$ws.Cells.Item(7, 14).Value = 5.395544
$ws.Cells.Item(7, 15).Value = 0.5830092258526149
$ws.Cells.Item(7, 16).Value = 0.5830092258526149
$ws.Cells.Item(7, 17).Value = 20.64231060351289
$ws.Cells.Item(7, 18).Value = 185.780795431616
$ws.Cells.Item(7, 19).Value = 0.4077568571155633
$ws.Cells.Item(7, 20).Value = 0.4077568571155632

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "H2-M3"
$ws.Cells.Item(8, 3).Value = "Klrd1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6224386666666667
$ws.Cells.Item(8, 8).Value = 1.867316
$ws.Cells.Item(8, 9).Value = 0.03792958446082136
$ws.Cells.Item(8, 10).Value = 0.03792958446082135
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 1.286367333333333
$ws.Cells.Item(8, 14).Value = 3.859102
$ws.Cells.Item(8, 15).Value = 0.416990774147385
$ws.Cells.Item(8, 16).Value = 0.4169907741473849
$ws.Cells.Item(8, 17).Value = 0.8006847678035557
$ws.Cells.Item(8, 18).Value = 7.206162910232
$ws.Cells.Item(8, 19).Value = 0.01581628678740652
$ws.Cells.Item(8, 20).Value = 0.01581628678740652

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "H2-M3"
$ws.Cells.Item(9, 3).Value = "Klrd1"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.6224386666666667
$ws.Cells.Item(9, 8).Value = 1.867316
$ws.Cells.Item(9, 9).Value = 0.03792958446082136
$ws.Cells.Item(9, 10).Value = 0.03792958446082135
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.798514666666667
$ws.Cells.Item(9, 14).Value = 5.395544
$ws.Cells.Item(9, 15).Value = 0.5830092258526149
$ws.Cells.Item(9, 16).Value = 0.5830092258526149
$ws.Cells.Item(9, 17).Value = 1.119465071100445
$ws.Cells.Item(9, 18).Value = 10.075185639904
$ws.Cells.Item(9, 19).Value = 0.02211329767341483
$ws.Cells.Item(9, 20).Value = 0.02211329767341483

